$wb = $excel.ActiveWorkbook

# --- Information sheet: move the selection from C3 to B3 -----------------
# (the C3 formula `=D3-NOW()` is volatile and recalculates on save, so its
# cached <v> naturally tracks wall-clock time without any extra action)
$wsInfo = $wb.Worksheets.Item("Information")
$null = $wsInfo.Range("B3").Select()

# --- UnityIsActive sheet: update the Controls panel text ------------------
$ws = $wb.Worksheets.Item("UnityIsActive")

# Re-word the "Excel input is blocked" banner
$ws.Range("B2").Value = "Excel input is currently blocked - access challenges through the overworld!"

# Re-order/re-word the control hints:
#   Interact:        Space/Return
#   Select dialogue:  Arrow keys
#   Menu:            M
$ws.Range("B8").Value = "Select dialogue:"
$ws.Range("B7").Value = "Interact:"
$ws.Range("C7").Value = "Space/Return"
$ws.Range("B9").Value = "Menu:"
$ws.Range("C9").Value = "M"
$ws.Range("C8").Value = "Arrow keys"

# Widen column B so the new labels fit
$ws.Columns.Item(2).ColumnWidth = 16.88

# Restore UnityIsActive as the active sheet/tab after touching other sheets
$null = $ws.Activate()
